$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price/volume data rows 2, 3, 4 and 9 are rotated:
#   new row2 = old row3, new row3 = old row4, new row4 = old row9, new row9 = old row2
# Apply the resulting final values directly to columns D, J, K, L, M, P.

$ws.Range("D2").Value = 44839
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15600
$ws.Range("P2").Value = 1040

$ws.Range("D3").Value = 44750
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 19000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 19571
$ws.Range("P3").Value = 1305

$ws.Range("D4").Value = 45133
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 22000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 22000
$ws.Range("P4").Value = 1467

$ws.Range("D9").Value = 44749
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 17556
$ws.Range("P9").Value = 1170
